# Auto-update draw results: append the 2025-11-24 Pick 4 result as a new
# row (row 69) at the bottom of the results table, matching the existing
# "stored as text" convention used for every other row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 69

# Force the new row's cells to Text format *before* assigning values so
# the date-looking and number-looking strings are not auto-converted by
# Excel into a real date serial / numeric value (the sheet otherwise
# stores every value - including dates and draw "phase" numbers - as
# plain text, suppressing the resulting "number stored as text" warning
# via <ignoredErrors>).
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("B$row").NumberFormat = "@"
$ws.Range("C$row").NumberFormat = "@"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("E$row").NumberFormat = "@"

$ws.Range("A$row").Value = "2025-11-24"
$ws.Range("B$row").Value = "Pick 4"
$ws.Range("C$row").Value = "251124"
$ws.Range("D$row").Value = "9-2-0-5"
$ws.Range("E$row").Value = "2025-11-24T21:41:02.041+04:00"
